$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header text updates (Volume Number + report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  38"
$ws.Range("C9").Value = "Report Covering the Week  9/15/2025  Through  9/21/2025"

# --- Row 15 ---
$ws.Range("D15").Value = 1
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -60
$ws.Range("J15").Value = 24
$ws.Range("K15").Value = 16.666666666666
$ws.Range("L15").Value = 154.545454545455
$ws.Range("N15").Value = -15.151515151515

# --- Row 16 ---
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -12.5
$ws.Range("F16").Value = 27
$ws.Range("G16").Value = 32
$ws.Range("H16").Value = -15.625
$ws.Range("I16").Value = 241
$ws.Range("J16").Value = 301
$ws.Range("K16").Value = -19.933554817275
$ws.Range("L16").Value = -29.941860465116
$ws.Range("M16").Value = 117.117117117117
$ws.Range("N16").Value = -87.139807897545

# --- Row 17 ---
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -75
$ws.Range("F17").Value = 36
$ws.Range("G17").Value = 39
$ws.Range("H17").Value = -7.692307692307
$ws.Range("I17").Value = 414
$ws.Range("J17").Value = 398
$ws.Range("K17").Value = 4.020100502512
$ws.Range("L17").Value = 11.891891891891
$ws.Range("M17").Value = 195.714285714286
$ws.Range("N17").Value = -15.510204081632

# --- Row 18 ---
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 36
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = 71.428571428571
$ws.Range("I18").Value = 318
$ws.Range("J18").Value = 268
$ws.Range("K18").Value = 18.656716417910
$ws.Range("L18").Value = 3.921568627450
$ws.Range("M18").Value = 34.745762711864
$ws.Range("N18").Value = -83.667180277349

# --- Row 19 ---
$ws.Range("C19").Value = 32
$ws.Range("D19").Value = 32
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 125
$ws.Range("G19").Value = 142
$ws.Range("H19").Value = -11.971830985915
$ws.Range("I19").Value = 1239
$ws.Range("J19").Value = 1463
$ws.Range("K19").Value = -15.311004784689
$ws.Range("L19").Value = -26.945754716981
$ws.Range("M19").Value = -24.173806609547
$ws.Range("N19").Value = -82.035667681600

# --- Row 20 ---
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -80
$ws.Range("J20").Value = 42
$ws.Range("K20").Value = -61.904761904761
$ws.Range("L20").Value = -69.811320754717
$ws.Range("N20").Value = -93.798449612403

# --- Row 21 ---
$ws.Range("C21").Value = 49
$ws.Range("D21").Value = 61
$ws.Range("E21").Value = -19.672131147541
$ws.Range("F21").Value = 227
$ws.Range("G21").Value = 244
$ws.Range("H21").Value = -6.967213114754
$ws.Range("I21").Value = 2257
$ws.Range("J21").Value = 2500
$ws.Range("K21").Value = -9.72
$ws.Range("L21").Value = -18.871315600287
$ws.Range("M21").Value = 5.123428039124
$ws.Range("N21").Value = -80.387556482447

# --- Row 22 ---
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 300
$ws.Range("F22").Value = 12
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 144
$ws.Range("J22").Value = 122
$ws.Range("K22").Value = 18.032786885245
$ws.Range("L22").Value = -7.096774193548
$ws.Range("M22").Value = 38.461538461538

# --- Row 24 ---
$ws.Range("C24").Value = 69
$ws.Range("D24").Value = 106
$ws.Range("E24").Value = -34.905660377358
$ws.Range("F24").Value = 324
$ws.Range("G24").Value = 334
$ws.Range("H24").Value = -2.994011976047
$ws.Range("I24").Value = 2940
$ws.Range("J24").Value = 3174
$ws.Range("K24").Value = -7.372400756143
$ws.Range("L24").Value = -0.338983050847
$ws.Range("M24").Value = -13.402061855670

# --- Row 25 ---
$ws.Range("C25").Value = 55
$ws.Range("D25").Value = 90
$ws.Range("E25").Value = -38.888888888888
$ws.Range("F25").Value = 259
$ws.Range("G25").Value = 286
$ws.Range("H25").Value = -9.440559440559
$ws.Range("I25").Value = 2498
$ws.Range("J25").Value = 2800
$ws.Range("K25").Value = -10.785714285714
$ws.Range("L25").Value = -6.511976047904

# --- Row 26 ---
$ws.Range("C26").Value = 21
$ws.Range("D26").Value = 24
$ws.Range("E26").Value = -12.5
$ws.Range("F26").Value = 96
$ws.Range("G26").Value = 81
$ws.Range("H26").Value = 18.518518518518
$ws.Range("I26").Value = 817
$ws.Range("J26").Value = 759
$ws.Range("K26").Value = 7.641633728590
$ws.Range("L26").Value = 4.342273307790
$ws.Range("M26").Value = 91.784037558685

# --- Row 27 ---
$ws.Range("D27").Value = 1
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -60
$ws.Range("J27").Value = 29
$ws.Range("K27").Value = 10.344827586206
$ws.Range("L27").Value = 60

# --- Row 28 ---
$ws.Range("C28").Value = 7
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = 133.333333333333
$ws.Range("F28").Value = 46
$ws.Range("G28").Value = 18
$ws.Range("H28").Value = 155.555555555556
$ws.Range("I28").Value = 294
$ws.Range("J28").Value = 165
$ws.Range("K28").Value = 78.181818181818
$ws.Range("L28").Value = 81.481481481481

# --- Row 29 (was "n/a" text in D/E/G/H, now numeric) ---
# Pull numeric formatting from this row's existing numeric cells before
# writing values, so the cells switch out of the text style correctly.
$ws.Range("I29").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("K29").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = -100
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -100
$ws.Range("J29").Value = 6
$ws.Range("K29").Value = -33.333333333333

# --- Row 30 (was "n/a" text in D/E/G/H, now numeric) ---
$ws.Range("I30").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("K30").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = -100
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = -100
$ws.Range("J30").Value = 6
$ws.Range("K30").Value = -66.666666666666

# --- Row 31 (was "n/a" text in D/E, now numeric; G/H already numeric) ---
$ws.Range("I31").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("K31").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = -100
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = -50
$ws.Range("J31").Value = 20
$ws.Range("K31").Value = -35

# --- Row 33 (was numeric in G/H, now "n/a" text) ---
# H33's new text ("***.*") isn't numeric-looking, so a plain format-paste
# followed by the value assignment keeps it text without any extra style.
$ws.Range("E33").Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("H33").Value = "***.*"

# G33's new text ("0") *is* numeric-looking, so it needs to be forced to
# the Text number format before assignment, then the real target format
# (matching the sibling "n/a" cells) is pasted back over it.
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "0"
$ws.Range("F33").Copy()
$ws.Range("G33").PasteSpecial(-4122)
